$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 78.5
$ws.Range("I11").Value = 78.5
$ws.Range("K11").Value = 78.5
$ws.Range("M11").Value = 61.5

$ws.Range("H17").Value = 2138.575
$ws.Range("J17").Value = 2138.575
$ws.Range("L17").Value = 6415.724999999999
$ws.Range("N17").Value = -6751.724999999999

$ws.Range("H43").Value = 6750
$ws.Range("J43").Value = 6000
$ws.Range("L43").Value = 6000
$ws.Range("N43").Value = -6138

$ws.Range("H53").Value = 133.33333
$ws.Range("I53").Value = 61.75
$ws.Range("J53").Value = 276.5
$ws.Range("K53").Value = 61.75
$ws.Range("L53").Value = 276.5
$ws.Range("M53").Value = 575.25
$ws.Range("N53").Value = -1550.5

$ws.Range("H64").Value = 4942.857
$ws.Range("J64").Value = 4920
$ws.Range("L64").Value = 4920
$ws.Range("N64").Value = -5416

$ws.Range("H67").Value = 4942.857
$ws.Range("J67").Value = 4920
$ws.Range("L67").Value = 4920
$ws.Range("N67").Value = -6636

$ws.Range("H100").Value = 1679.4

$ws.Range("H132").Value = 2546.5715
$ws.Range("I132").Value = 2721.2778
$ws.Range("J132").Value = 1498.3334
$ws.Range("K132").Value = 8163.8334
$ws.Range("L132").Value = 4495.0002
$ws.Range("M132").Value = -5633.8334
$ws.Range("N132").Value = -9555.0002

$ws.Range("H137").Value = 1908.4445
$ws.Range("I137").Value = 1539.1428
$ws.Range("J137").Value = 3201
$ws.Range("K137").Value = 4617.428400000001
$ws.Range("L137").Value = 9603
$ws.Range("M137").Value = -2067.428400000001
$ws.Range("N137").Value = -14703

$ws.Range("H138").Value = 3572.6047
$ws.Range("J138").Value = 2602.3242
$ws.Range("L138").Value = 7806.9726
$ws.Range("N138").Value = -18086.9726

$ws.Range("H141").Value = 4216.1113
$ws.Range("I141").Value = 3191.4285
$ws.Range("K141").Value = 9574.2855
$ws.Range("M141").Value = -4394.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 10.833333
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = 107
$ws.Range("N4").Value = -252

$ws.Range("H5").Value = 173.77777
$ws.Range("I5").Value = 109.14286
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 109.14286
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = 2.857140000000001
$ws.Range("N5").Value = -624

$ws.Range("H32").Value = 5601.806
$ws.Range("I32").Value = 3764.6545
$ws.Range("J32").Value = 14022.083
$ws.Range("K32").Value = 3764.6545
$ws.Range("L32").Value = 14022.083
$ws.Range("M32").Value = -3477.6545
$ws.Range("N32").Value = -14596.083

$ws.Range("H74").Value = 1886.9
$ws.Range("I74").Value = 1141.1765
$ws.Range("K74").Value = 1141.1765
$ws.Range("M74").Value = -267.1765

$ws.Range("H77").Value = 1886.9
$ws.Range("I77").Value = 1141.1765
$ws.Range("K77").Value = 5705.8825
$ws.Range("M77").Value = -1337.8825

$ws.Range("H101").Value = 27499
$ws.Range("J101").Value = 27499
$ws.Range("L101").Value = 27499
$ws.Range("N101").Value = -33989

$ws.Range("H106").Value = 95000
$ws.Range("J106").Value = 95000
$ws.Range("L106").Value = 95000
$ws.Range("N106").Value = -97524

$ws.Range("H132").Value = 2387.913
$ws.Range("I132").Value = 2372.2104
$ws.Range("K132").Value = 7116.6312
$ws.Range("M132").Value = -4586.6312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 173.77777
$ws.Range("I4").Value = 109.14286
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 109.14286
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = 5.857140000000001
$ws.Range("N4").Value = -630

$ws.Range("H22").Value = 663.7273
$ws.Range("I22").Value = 663.7273
$ws.Range("K22").Value = 663.7273
$ws.Range("M22").Value = -490.7273

$ws.Range("H141").Value = 49998.5
$ws.Range("J141").Value = 49998
$ws.Range("L141").Value = 49998
$ws.Range("N141").Value = -60358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 179.65715
$ws.Range("I7").Value = 130.34483
$ws.Range("K7").Value = 130.34483
$ws.Range("M7").Value = -17.34483

$ws.Range("H22").Value = 599
$ws.Range("I22").Value = 599
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 599
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -249

$ws.Range("H41").Value = 1054.5
$ws.Range("I41").Value = 1054.5
$ws.Range("K41").Value = 1054.5
$ws.Range("M41").Value = -626.5

$ws.Range("H99").Value = 12445.762
$ws.Range("I99").Value = 8084.8
$ws.Range("J99").Value = 16410.273
$ws.Range("K99").Value = 8084.8
$ws.Range("L99").Value = 16410.273
$ws.Range("M99").Value = -6586.8
$ws.Range("N99").Value = -19406.273

$ws.Range("H107").Value = 743.1667
$ws.Range("I107").Value = 702.25
$ws.Range("J107").Value = 825
$ws.Range("K107").Value = 702.25
$ws.Range("L107").Value = 825
$ws.Range("M107").Value = 1217.75
$ws.Range("N107").Value = -4665

$ws.Range("H126").Value = 12445.762
$ws.Range("I126").Value = 8084.8
$ws.Range("J126").Value = 16410.273
$ws.Range("K126").Value = 24254.4
$ws.Range("L126").Value = 49230.819
$ws.Range("M126").Value = -21784.4
$ws.Range("N126").Value = -54170.819

$ws.Range("H132").Value = 2600.3684
$ws.Range("I132").Value = 2600.3684
$ws.Range("K132").Value = 7801.1052
$ws.Range("M132").Value = -5271.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66682.92999999999
$ws.Range("I2").Value = 100008.2
$ws.Range("J2").Value = 32.4
$ws.Range("K2").Value = 600049.2
$ws.Range("L2").Value = 194.4
$ws.Range("M2").Value = -599936.2
$ws.Range("N2").Value = -420.4

$ws.Range("H29").Value = 500303.5
$ws.Range("I29").Value = 1000047.5
$ws.Range("J29").Value = 559.5
$ws.Range("K29").Value = 3000142.5
$ws.Range("L29").Value = 1678.5
$ws.Range("M29").Value = -2999865.5
$ws.Range("N29").Value = -2232.5

$ws.Range("H34").Value = 1865.091
$ws.Range("J34").Value = 3200
$ws.Range("L34").Value = 9600
$ws.Range("N34").Value = -9768

$ws.Range("H46").Value = 3333966.8
$ws.Range("J46").Value = 3333966.8
$ws.Range("L46").Value = 10001900.4
$ws.Range("N46").Value = -10002082.4

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null

$ws.Range("H129").Value = 901.4286
$ws.Range("J129").Value = 700
$ws.Range("L129").Value = 2100
$ws.Range("N129").Value = -12100

$ws.Range("H134").Value = 2918.1428
$ws.Range("I134").Value = 2918.1428
$ws.Range("K134").Value = 8754.428400000001
$ws.Range("M134").Value = -3684.428400000001

$ws.Range("H136").Value = 10059.333
$ws.Range("I136").Value = 5089
$ws.Range("K136").Value = 15267
$ws.Range("M136").Value = -10167

$ws.Range("H137").Value = 4849.6875
$ws.Range("J137").Value = 4784.7
$ws.Range("L137").Value = 14354.1
$ws.Range("N137").Value = -24554.1

$ws.Range("H138").Value = 7399.6
$ws.Range("I138").Value = 5333
$ws.Range("K138").Value = 15999
$ws.Range("M138").Value = -10859

$ws.Range("H139").Value = 8999.666999999999
$ws.Range("I139").Value = 6000
$ws.Range("J139").Value = 10499.5
$ws.Range("K139").Value = 18000
$ws.Range("L139").Value = 31498.5
$ws.Range("M139").Value = -12860
$ws.Range("N139").Value = -41778.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4826.2856
$ws.Range("I132").Value = 3199
$ws.Range("J132").Value = 5477.2
$ws.Range("K132").Value = 9597
$ws.Range("L132").Value = 16431.6
$ws.Range("M132").Value = -7067
$ws.Range("N132").Value = -21491.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5118.9165
$ws.Range("J22").Value = 5947.8
$ws.Range("L22").Value = 5947.8
$ws.Range("N22").Value = -6537.8

$ws.Range("H27").Value = 5118.9165
$ws.Range("J27").Value = 5947.8
$ws.Range("L27").Value = 5947.8
$ws.Range("N27").Value = -6161.8

$ws.Range("H74").Value = 37500
$ws.Range("I74").Value = 25000
$ws.Range("K74").Value = 25000
$ws.Range("M74").Value = -24002

$ws.Range("H77").Value = 37500
$ws.Range("I77").Value = 25000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70008

$ws.Range("H106").Value = 21249.75
$ws.Range("J106").Value = 21249.75
$ws.Range("L106").Value = 21249.75
$ws.Range("N106").Value = -23773.75

$ws.Range("H136").Value = 2868.7896
$ws.Range("I136").Value = 2772.6667
$ws.Range("K136").Value = 8318.000100000001
$ws.Range("M136").Value = -5768.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3212
$ws.Range("I122").Value = 3212
$ws.Range("K122").Value = 9636
$ws.Range("M122").Value = -7186

$ws.Range("H132").Value = 1994.25
$ws.Range("I132").Value = 1992.3334
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5977.0002
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3447.0002
$ws.Range("N132").Value = -11060
